$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (shared between both sheets' B1 cell) to add "-1st" suffix
$wsInput.Range("B1").Value = "2601-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-VALIDATE-RANGE-1st"
$wsOutput.Range("B1").Value = "2601-RBI-EI-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-VALIDATE-RANGE-1st"

# Change shortname cell from numeric 2601 to text "260d"
$wsInput.Range("B2").Value = "260d"

# Make ProductLoanInput the selected/active sheet instead of ProductLoanOutput
$wsInput.Select()
